# Automatische test-sync: 2025-08-13 21:45:50
#
# Appends the new "Demo inplannen" log entry (row 13) to the "Logs" sheet,
# and bumps the corresponding tally on the "Dashboard" sheet from 11 to 12.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 13

$logs.Cells.Item($newRow, 1).Value = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-13 21:45:14"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 12

# Extend the conditional-formatting ranges (D, G, H, I, J) from row 12 to row 13
# now that a new data row has been appended.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`12")
    $newRange = $logs.Range("$col`2:$col`13")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
